$p = $ppt.ActivePresentation

# --- 1. Update the cached "datetimeFigureOut" date field text from
#        7/7/2021 to 7/14/2021 everywhere it appears: on the slide
#        master's Date Placeholder and on each of the 11 slide layouts'
#        Date Placeholder shapes. ---

$master = $p.SlideMaster

function Set-DatePlaceholderText {
    param($shapes)
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "7/14/2021"
        }
    }
}

Set-DatePlaceholderText $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lay = $layouts.Item($i)
    Set-DatePlaceholderText $lay.Shapes
}

# --- 2. Swap the M2/M3 text-box labels on slide 1 and widen both boxes
#        from 498855 EMU to 503664 EMU (39.27992pt -> 39.6586pt). ---

$s1 = $p.Slides.Item(1)
$shpM2 = $null   # "TextBox 141", currently reads "M2"
$shpM3 = $null   # "TextBox 142", currently reads "M3"
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 141") { $shpM2 = $shp }
    elseif ($shp.Name -eq "TextBox 142") { $shpM3 = $shp }
}

$shpM2.TextFrame.TextRange.Text = "M3"
$shpM3.TextFrame.TextRange.Text = "M2"

$shpM2.Width = 39.6586
$shpM3.Width = 39.6586
